# Add files via upload
# Append the "17.10 / 20.10 / 23.10" task blocks to the Zeitmanagement sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

# ---------------------------------------------------------------------------
# Block 1: row 33-35  (Samstag, 17.10.2020 -> serial 44121)
# ---------------------------------------------------------------------------
$ws.Range("A14").Copy()
$ws.Range("A33").PasteSpecial($xlPasteFormats)
$ws.Range("A33").Value = 44121

$ws.Range("B1:E1").Copy()
$ws.Range("B33:E33").PasteSpecial($xlPasteFormats)
$ws.Range("B33:E33").Merge()
$ws.Range("B33").Value = "Recherche Google (verschiedene Suchbegriffe benutzt)"

$ws.Range("F14:G14").Copy()
$ws.Range("F33:G33").PasteSpecial($xlPasteFormats)
$ws.Range("F33").Value = 0.5
$ws.Range("G33").Value = 0.5625

$ws.Range("I14").Copy()
$ws.Range("I33").PasteSpecial($xlPasteFormats)
$ws.Range("I33").Value = 90

$ws.Range("B1:E1").Copy()
$ws.Range("B34:E34").PasteSpecial($xlPasteFormats)
$ws.Range("B34:E34").Merge()
$ws.Range("B34").Value = "Alte Links durchgegangen "

$ws.Range("F14:G14").Copy()
$ws.Range("F34:G34").PasteSpecial($xlPasteFormats)
$ws.Range("F34").Value = 0.58333333333333337
$ws.Range("G34").Value = 0.66666666666666663

$ws.Range("I14").Copy()
$ws.Range("I34").PasteSpecial($xlPasteFormats)
$ws.Range("I34").Value = 120

$ws.Range("J16").Copy()
$ws.Range("J35").PasteSpecial($xlPasteFormats)
$ws.Range("J35").Value = 210

$ws.Range("H23").Copy()
$ws.Range("K35").PasteSpecial($xlPasteFormats)
$ws.Range("K35").Value = "3h 30min"

# ---------------------------------------------------------------------------
# Block 2: row 37-38  (Dienstag, 20.10.2020 -> serial 44124)
# ---------------------------------------------------------------------------
$ws.Range("A14").Copy()
$ws.Range("A37").PasteSpecial($xlPasteFormats)
$ws.Range("A37").Value = 44124

$ws.Range("B1:E1").Copy()
$ws.Range("B37:E37").PasteSpecial($xlPasteFormats)
$ws.Range("B37:E37").Merge()
$ws.Range("B37").Value = "GUI mit Exit Button verbessert"

$ws.Range("F14:G14").Copy()
$ws.Range("F37:G37").PasteSpecial($xlPasteFormats)
$ws.Range("F37").Value = 0.5625
$ws.Range("G37").Value = 0.56944444444444442

$ws.Range("I14").Copy()
$ws.Range("I37").PasteSpecial($xlPasteFormats)
$ws.Range("I37").Value = 10

$ws.Range("B1:E1").Copy()
$ws.Range("B38:E38").PasteSpecial($xlPasteFormats)
$ws.Range("B38:E38").Merge()
$ws.Range("B38").Value = "veruscht einen Algorithmus zu programmieren [abgebrochen]"

$ws.Range("F14:G14").Copy()
$ws.Range("F38:G38").PasteSpecial($xlPasteFormats)
$ws.Range("F38").Value = 0.56944444444444442
$ws.Range("G38").Value = 0.625

$ws.Range("I14").Copy()
$ws.Range("I38").PasteSpecial($xlPasteFormats)
$ws.Range("I38").Value = 80

# ---------------------------------------------------------------------------
# Block 3: row 39-46  (Freitag, 23.10.2020 -> serial 44127)
# ---------------------------------------------------------------------------
$ws.Range("A14").Copy()
$ws.Range("A39").PasteSpecial($xlPasteFormats)
$ws.Range("A39").Value = 44127

$ws.Range("B1:E1").Copy()
$ws.Range("B39:E39").PasteSpecial($xlPasteFormats)
$ws.Range("B39:E39").Merge()
$ws.Range("B39").Value = "OpenCV Download und Installation"

$ws.Range("F14:G14").Copy()
$ws.Range("F39:G39").PasteSpecial($xlPasteFormats)
$ws.Range("F39").Value = 0.63541666666666663
$ws.Range("G39").Value = 0.64583333333333337

$ws.Range("I14").Copy()
$ws.Range("I39").PasteSpecial($xlPasteFormats)
$ws.Range("I39").Value = 15

$ws.Range("B1:E1").Copy()
$ws.Range("B40:E40").PasteSpecial($xlPasteFormats)
$ws.Range("B40:E40").Merge()
$ws.Range("B40").Value = "OpenCV mit Netbeans verbunden "

$ws.Range("F14:G14").Copy()
$ws.Range("F40:G40").PasteSpecial($xlPasteFormats)
$ws.Range("F40").Value = 0.64583333333333337
$ws.Range("G40").Value = 0.65972222222222221

$ws.Range("I14").Copy()
$ws.Range("I40").PasteSpecial($xlPasteFormats)
$ws.Range("I40").Value = 20

$ws.Range("B1:E1").Copy()
$ws.Range("B41:E41").PasteSpecial($xlPasteFormats)
$ws.Range("B41:E41").Merge()
$ws.Range("B41").Value = "Kurzbericht erstellt und abgegeben"

$ws.Range("F14:G14").Copy()
$ws.Range("F41:G41").PasteSpecial($xlPasteFormats)
$ws.Range("F41").Value = 0.65972222222222221
$ws.Range("G41").Value = 0.66319444444444442

$ws.Range("I14").Copy()
$ws.Range("I41").PasteSpecial($xlPasteFormats)
$ws.Range("I41").Value = 5

$ws.Range("B1:E1").Copy()
$ws.Range("B42:E42").PasteSpecial($xlPasteFormats)
$ws.Range("B42:E42").Merge()
$ws.Range("B42").Value = "Stundenaufzeichnung aktualiesiert und abgegeben"

$ws.Range("F14:G14").Copy()
$ws.Range("F42:G42").PasteSpecial($xlPasteFormats)
$ws.Range("F42").Value = 0.66319444444444442
$ws.Range("G42").Value = 0.66666666666666663

$ws.Range("I14").Copy()
$ws.Range("I42").PasteSpecial($xlPasteFormats)
$ws.Range("I42").Value = 5

$ws.Range("B1:E1").Copy()
$ws.Range("B43:E43").PasteSpecial($xlPasteFormats)
$ws.Range("B43:E43").Merge()
$ws.Range("B43").Value = "Auf GitHub gepusht"

$ws.Range("F14:G14").Copy()
$ws.Range("F43:G43").PasteSpecial($xlPasteFormats)
$ws.Range("F43").Value = 0.66666666666666663
$ws.Range("G43").Value = 0.67361111111111116

$ws.Range("I14").Copy()
$ws.Range("I43").PasteSpecial($xlPasteFormats)
$ws.Range("I43").Value = 10

$ws.Range("B1:E1").Copy()
$ws.Range("B44:E44").PasteSpecial($xlPasteFormats)
$ws.Range("B44:E44").Merge()
$ws.Range("B44").Value = "Formular Pflichtpraktikum aktualiesiert "

$ws.Range("F14:G14").Copy()
$ws.Range("F44:G44").PasteSpecial($xlPasteFormats)
$ws.Range("F44").Value = 0.67361111111111116
$ws.Range("G44").Value = 0.67708333333333337

$ws.Range("I14").Copy()
$ws.Range("I44").PasteSpecial($xlPasteFormats)
$ws.Range("I44").Value = 5

$ws.Range("B1:E1").Copy()
$ws.Range("B45:E45").PasteSpecial($xlPasteFormats)
$ws.Range("B45:E45").Merge()
$ws.Range("B45").Value = "Eigenes Protokoll laufend aktualisiert"

$ws.Range("F14:G14").Copy()
$ws.Range("F45:G45").PasteSpecial($xlPasteFormats)
$ws.Range("F45").Value = 0.67708333333333337
$ws.Range("G45").Value = 0.68402777777777779

$ws.Range("I14").Copy()
$ws.Range("I45").PasteSpecial($xlPasteFormats)
$ws.Range("I45").Value = 10

$ws.Range("J16").Copy()
$ws.Range("J46").PasteSpecial($xlPasteFormats)
$ws.Range("J46").Value = 150

$ws.Range("K21").Copy()
$ws.Range("K46").PasteSpecial($xlPasteFormats)
$ws.Range("K46").Value = "2h 30min"

# Clear the clipboard marquee / selection artifact left by the copy chain and
# park the selection roughly where the author left it before saving.
$ws.Range("K49").Select()

Write-Host "done"
